$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp note in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 23:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1124676
$ws.Range("C4").Value = 29653
$ws.Range("D4").Value = 158993
$ws.Range("E4").Value = 900173
$ws.Range("F4").Value = 16474
$ws.Range("G4").Value = 1654
$ws.Range("H4").Value = 65510

# Row 9 - Alemania
$ws.Range("B9").Value = 163936
$ws.Range("C9").Value = 927
$ws.Range("E9").Value = 30328
$ws.Range("G9").Value = 85
$ws.Range("H9").Value = 6708

# Row 104 - Burkina Faso
$ws.Range("B104").Value = 649
$ws.Range("C104").Value = 4
$ws.Range("D104").Value = 517
$ws.Range("E104").Value = 88
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 44

# Row 147 - Cabo Verde
$ws.Range("D147").Value = 18
$ws.Range("E147").Value = 103
